$d = $word.ActiveDocument
$d.Content.Find.Execute("пар между", $true, $false, $false, $false, $false,
                         $true, 1, $false, "между", 2)
